$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string (A1): 14:30 -> 15:47
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 15:47"

# Reorder "Papua Nueva Guinea": it moves from after "Islas Turcas y Caicos" to
# right after "Bermudas" (row 181), with brand-new statistics. Implemented as
# insert-new-row-182 (pushing San Martin..Islas Turcas down by one) followed
# by deleting the old Papua Nueva Guinea row (which is now at row 189).
$ws.Rows.Item(182).Insert()
$ws.Range("A182").Value = "Papua Nueva Guinea"
$ws.Range("B182").Value = 153
$ws.Range("C182").Value = 39
$ws.Range("D182").Value = 44
$ws.Range("E182").Value = 107
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 2
$ws.Rows.Item(189).Delete()

# Updated COVID statistics for several countries (Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes).

# Estados Unidos (row 4)
$ws.Range("B4").Value = 4920095
$ws.Range("C4").Value = 1675
$ws.Range("D4").Value = 2482909
$ws.Range("E4").Value = 2276821
$ws.Range("G4").Value = 75
$ws.Range("H4").Value = 160365

# India (row 6)
$ws.Range("B6").Value = 1920719
$ws.Range("C6").Value = 14106
$ws.Range("D6").Value = 1291757
$ws.Range("E6").Value = 589012
$ws.Range("G6").Value = 130
$ws.Range("H6").Value = 39950

# Arabia Saudita (row 16)
$ws.Range("B16").Value = 282824
$ws.Range("C16").Value = 1389
$ws.Range("D16").Value = 245314
$ws.Range("E16").Value = 34490
$ws.Range("G16").Value = 36
$ws.Range("H16").Value = 3020

# Alemania (row 22)
$ws.Range("B22").Value = 213309
$ws.Range("C22").Value = 229
$ws.Range("E22").Value = 9369
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = 9240

# Irak (row 24)
$ws.Range("B24").Value = 137556
$ws.Range("C24").Value = 2834
$ws.Range("D24").Value = 98442
$ws.Range("E24").Value = 34020
$ws.Range("G24").Value = 77
$ws.Range("H24").Value = 5094

# Suecia (row 34)
$ws.Range("B34").Value = 81540
$ws.Range("C34").Value = 66
$ws.Range("G34").Value = 5
$ws.Range("H34").Value = 5760

# Portugal (row 48)
$ws.Range("B48").Value = 51848
$ws.Range("C48").Value = 167
$ws.Range("D48").Value = 37565
$ws.Range("E48").Value = 12543
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 1740

# Uzbekistan (row 61)
$ws.Range("B61").Value = 27554
$ws.Range("C61").Value = 507
$ws.Range("D61").Value = 18682
$ws.Range("E61").Value = 8703
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 169

# Serbia (row 63)
$ws.Range("B63").Value = 27033
$ws.Range("C63").Value = 295
$ws.Range("E63").Value = 12372
$ws.Range("G63").Value = 9
$ws.Range("H63").Value = 614

# Estado de Palestina (row 80)
$ws.Range("E80").Value = 6358
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = 89

# Republica de Macedonia (row 84)
$ws.Range("B84").Value = 11289
$ws.Range("C84").Value = 87
$ws.Range("D84").Value = 7221
$ws.Range("E84").Value = 3557
$ws.Range("G84").Value = 6
$ws.Range("H84").Value = 511

# Noruega (row 86)
$ws.Range("B86").Value = 9388
$ws.Range("C86").Value = 26
$ws.Range("E86").Value = 380
